# Updated AddChargeDialogSlotFunctions to set charge.type (i.e. moving).
# Rename the "Offense Type" values in column D:
#   "Moving Traffic"     -> "Moving"
#   "Non-moving Traffic" -> "Non-moving"
# ("Criminal" values are left untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 34

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)   # column D
    $val = $cell.Value()
    if ($val -eq "Moving Traffic") {
        $cell.Value = "Moving"
    } elseif ($val -eq "Non-moving Traffic") {
        $cell.Value = "Non-moving"
    }
}
